# Fruta / hortaliza, semanal
# A new weekly price record is inserted at row 180 (pushing all the
# existing data rows down by one), adding the latest observation for
# "Poroto granado" at Femacal de La Calera (Coquimbo).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 180; everything below shifts down.
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row with the latest weekly observation.
$ws.Cells.Item(180, 1).Value  = 3
$ws.Cells.Item(180, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(180, 3).Value  = "Coquimbo"
$ws.Cells.Item(180, 4).Value  = 44932
$ws.Cells.Item(180, 5).Value  = 5
$ws.Cells.Item(180, 6).Value  = 100112030
$ws.Cells.Item(180, 7).Value  = "Poroto granado"
$ws.Cells.Item(180, 8).Value  = "Sin especificar"
$ws.Cells.Item(180, 9).Value  = "Primera"
$ws.Cells.Item(180, 10).Value = 73
$ws.Cells.Item(180, 11).Value = 39000
$ws.Cells.Item(180, 12).Value = 40000
$ws.Cells.Item(180, 13).Value = 39479
$ws.Cells.Item(180, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(180, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(180, 16).Value = 1579
$ws.Cells.Item(180, 17).Value = 25
$ws.Cells.Item(180, 18).Value = "Hortaliza"
